# Pharma_Society_Report.xlsx - membership count refresh
#
# The worksheet tab is renamed from "Report" to "Sheet1", and the
# membership-count column (B) for each society is updated to the latest
# figures. The underlying society rows/order (FLASCO, GASCO, IOS, IOWA,
# MOASC) and all other columns are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Report" -> "Sheet1"
$ws.Name = "Sheet1"

# Refresh membership counts (column B) for each society row
$ws.Range("B2").Value = 850   # FLASCO (Florida Society of Clinical Oncology)
$ws.Range("B3").Value = 450   # GASCO (Georgia Society of Clinical Oncology)
$ws.Range("B4").Value = 250   # IOS (Indiana Oncology Society)
$ws.Range("B5").Value = 137   # IOWA Oncology Society
$ws.Range("B6").Value = 650   # MOASC (Medical Oncology Association of Southern California)

# Reset the active selection back to the top-left cell
$ws.Range("A1").Select()
